# Fix testdata: columns I and J ("Kortudsteder" / "Hændelsestype") were
# swapped in the source data - the event-type value ("Køb") was stored in
# the card-issuer column (J) instead of the event-type column (I).
# Swap the two columns' contents back (header row 17 + data rows 18-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 17; $row -le 23; $row++) {
    $cellI = $ws.Cells.Item($row, 9)   # column I
    $cellJ = $ws.Cells.Item($row, 10)  # column J

    $valI = $cellI.Value()
    $valJ = $cellJ.Value()

    $cellI.Value = $valJ
    $cellJ.Value = $valI
}

# Column width tweaks (OOXML stored widths, expressed here as Excel
# COM "number of characters" ColumnWidth - stored = ColumnWidth + 5/6,
# quantized to 1/6 character-width steps).
$ws.Columns.Item(5).ColumnWidth = (6.42 - 5/6)
$ws.Columns.Item(6).ColumnWidth = (11.25 - 5/6)
$ws.Columns.Item(7).ColumnWidth = (13.06 - 5/6)
$ws.Columns.Item(9).ColumnWidth = (12.96 - 5/6)
$ws.Columns.Item(10).ColumnWidth = (19.95 - 5/6)
$ws.Columns.Item(11).ColumnWidth = (12.96 - 5/6)

# Move the active selection from H32 to H9.
$ws.Range("H9").Select()
